# Raw and Clean Data from SSA for June 9th (date serial 43991)
# Adds the June 9, 2020 row to each of the daily tracking sheets and
# fills in the corresponding "control_obs" comparison column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) out_vars — new row 10 (2020-06-09 daily summary)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("out_vars")

# Column A carries the date-formatted style already used by A2:A9 — copy
# that formatting down before overwriting the value so the new cell picks
# up the same number format / font as the rest of the column.
$ws1.Range("A9").Copy($ws1.Range("A10"))
$ws1.Range("A10").Value = 43991

$ws1.Range("B10").Value = 124301
$ws1.Range("B10").WrapText = $true
$ws1.Range("C10").Value = 182077
$ws1.Range("C10").WrapText = $true
$ws1.Range("D10").Value = 50677
$ws1.Range("D10").WrapText = $true
$ws1.Range("E10").Value = 14649
$ws1.Range("E10").WrapText = $true
$ws1.Range("F10").Value = 33.206490695971873
$ws1.Range("F10").WrapText = $true
$ws1.Range("G10").Value = 41276
$ws1.Range("G10").WrapText = $true
$ws1.Range("H10").Value = 3826
$ws1.Range("H10").WrapText = $true
$ws1.Range("I10").Value = 3982
$ws1.Range("I10").WrapText = $true
$ws1.Range("J10").Value = 357055
$ws1.Range("J10").WrapText = $true

$ws1.Activate()
$ws1.Range("D26").Select()

# ---------------------------------------------------------------------
# 2) dates_dx — fill the already-present (blank) row 10
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dates_dx")
$ws2.Range("A10").Value = 43991
$ws2.Range("B10").Value = 0
$ws2.Range("C10").Value = 1
$ws2.Range("D10").Value = 1
$ws2.Range("E10").Value = 1
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 0
$ws2.Range("H10").Value = 0
$ws2.Range("I10").Value = 4

$ws2.Activate()
$ws2.Range("D36").Select()

# ---------------------------------------------------------------------
# 3) dates_sx — fill row 10 (only A10 existed before, B10:L10 are new)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dates_sx")
$ws3.Range("A10").Value = 43991
$ws3.Range("B10").Value = 0
$ws3.Range("C10").Value = 1
$ws3.Range("D10").Value = 0
$ws3.Range("E10").Value = 1
$ws3.Range("F10").Value = 1
$ws3.Range("G10").Value = 1
$ws3.Range("H10").Value = 0
$ws3.Range("I10").Value = 1
$ws3.Range("J10").Value = 1
$ws3.Range("K10").Value = 0
$ws3.Range("L10").Value = 0

$ws3.Activate()
$ws3.Range("M10").Select()

# ---------------------------------------------------------------------
# 4) dates_deaths — new row 10 plus 23 placeholder date rows (11:33)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dates_deaths")
$ws4.Range("A10").Value = 43991
$ws4.Range("B10").Value = 1
$ws4.Range("C10").Value = 0
$ws4.Range("D10").Value = 2
$ws4.Range("E10").Value = 1
$ws4.Range("F10").Value = 1
$ws4.Range("G10").Value = 2
$ws4.Range("H10").Value = 2

# Rows 11:33 only carry column A with the date number format (mm-dd-yy ==
# builtin numFmtId 14), no values yet.
$ws4.Range("A11").NumberFormat = "mm-dd-yy"
$ws4.Range("A11").Copy($ws4.Range("A12:A33"))

$ws4.Activate()
$ws4.Range("I10").Select()

# ---------------------------------------------------------------------
# 5) control_obs — new "2020-06-09" column (J)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("control_obs")
$ws5.Range("J1").Value = 43991
$ws5.Range("J2").Value = 3305
$ws5.Range("J3").Value = 3120
$ws5.Range("J4").Value = 3120
$ws5.Range("J5").Value = 3120
$ws5.Range("J6").Value = 3120
$ws5.Range("J7").Value = 2395
$ws5.Range("J8").Value = 4942
$ws5.Range("J10").Value = 149
$ws5.Range("J11").Value = 149
$ws5.Range("J12").Value = 149
$ws5.Range("J13").Value = 149
$ws5.Range("J14").Value = 149
$ws5.Range("J15").Value = 126
$ws5.Range("J16").Value = 161
$ws5.Range("J18").Value = 799

$ws5.Activate()
$ws5.Range("J18").Select()

# Leave the workbook focused back on out_vars (the tab that was active
# both before and after this edit).
$ws1.Activate()
